# Insert a new data row at sheet row 36 (pushing existing rows 36:132 down to 37:133)
# and populate it with the new "Arveja Verde" record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(36).Insert()

$ws.Range("A36").Value = 10
$ws.Range("B36").Value = "Vega Modelo de Temuco"
$ws.Range("C36").Value = "La Araucanía"
$ws.Range("D36").Value = 45260
$ws.Range("E36").Value = 9
$ws.Range("F36").Value = 100112022
$ws.Range("G36").Value = "Arveja Verde"
$ws.Range("H36").Value = "Sin especificar"
$ws.Range("I36").Value = "Primera"
$ws.Range("J36").Value = 60
$ws.Range("K36").Value = 28000
$ws.Range("L36").Value = 28000
$ws.Range("M36").Value = 28000
$ws.Range("N36").Value = "$/saco 25 kilos"
$ws.Range("O36").Value = "Región del Maule"
$ws.Range("P36").Value = 1120
$ws.Range("Q36").Value = 25
$ws.Range("R36").Value = "Hortaliza"
